# Item-upload-template.xlsx: drop the Supplier / Supplier Code / Defect Type
# columns from the header row. The "Fields marked * are required." note
# (previously column F) slides left into column C as columns C:E are removed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("C:E").Delete()

# The tall row height was only needed for the wrapped "Defect Type" header;
# with that column gone, let the row go back to its natural height.
$ws.Rows(1).AutoFit()

# Restore the (now different) selected cell that was active when the sheet
# was last saved.
$ws.Range("C6").Select()
